# FDR Calculation sheet automation.
# Adds a "Actual Result" -> Pass/Fail automation column (H) that reads
# "Passed" (green-ish highlight) for the rows whose check already succeeded
# (column G = "Pass") and "Failed" (highlight) for the row whose check
# failed (column G = "Fail"), completing the automation of the FDR sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5 passed their validation (column G = "Pass") -> mark column H "Passed"
$passedRange = $ws.Range("H2:H5")
$passedRange.Value = "Passed"
$passedRange.Interior.ColorIndex = 17
$passedRange.Borders.LineStyle = 1

# Row 6 failed its validation (column G = "Fail") -> mark column H "Failed"
$failedRange = $ws.Range("H6")
$failedRange.Value = "Failed"
$failedRange.Interior.ColorIndex = 10
$failedRange.Borders.LineStyle = 1

# Leave the selection where the author left it when they finished editing.
[void]$ws.Range("K9").Select()
